$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.377.84"
$ws.Range("E2").Value = "  +0.61%  "

$ws.Range("D3").Value = "1.863.13"
$ws.Range("E3").Value = "  +0.41%  "

$ws.Range("E4").Value = "  -0.26%  "

$ws.Range("D5").Value = "'246.38"
$ws.Range("E5").Value = "  +2.17%  "

$ws.Range("D6").Value = "'0.7045"
$ws.Range("E6").Value = "  +0.74%  "

$ws.Range("D7").Value = "'0.9978"
$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("D8").Value = "'0.07767"
$ws.Range("E8").Value = "  -0.14%  "

$ws.Range("D9").Value = "'0.3084"
$ws.Range("E9").Value = "  +0.35%  "

$ws.Range("D10").Value = "'23.84"
$ws.Range("E10").Value = "  +0.55%  "

$ws.Range("D11").Value = "'0.07830"
$ws.Range("E11").Value = "  +0.30%  "

$ws.Range("D12").Value = "'5.173"
$ws.Range("E12").Value = "  +1.29%  "

$ws.Range("D13").Value = "'93.45"
$ws.Range("E13").Value = "  +1.50%  "

$ws.Range("D14").Value = "1.850.71"
$ws.Range("E14").Value = "  -0.39%  "

$ws.Range("D15").Value = "'0.6950"
$ws.Range("E15").Value = "  +1.14%  "

$ws.Range("D16").Value = "'6.649"
$ws.Range("E16").Value = "  +1.52%  "

$ws.Range("D17").Value = "'0.000008371"
$ws.Range("E17").Value = "  -0.97%  "

$ws.Range("D18").Value = "29.343.44"
$ws.Range("E18").Value = "  +0.49%  "

$ws.Range("D19").Value = "'243.92"
$ws.Range("E19").Value = "  -1.65%  "

$ws.Range("D20").Value = "2.095.91"
$ws.Range("E20").Value = "  -0.66%  "

$ws.Range("D21").Value = "'12.81"
$ws.Range("E21").Value = "  -0.06%  "

$ws.Range("E22").Value = "  -0.20%  "

$ws.Range("D23").Value = "'7.591"
$ws.Range("E23").Value = "  +0.66%  "

$ws.Range("D24").Value = "'0.9976"
$ws.Range("E24").Value = "  -0.27%  "

$ws.Range("D25").Value = "'0.1522"
$ws.Range("E25").Value = "  +0.94%  "

$ws.Range("D26").Value = "'8.921"
$ws.Range("E26").Value = "  +0.71%  "

$ws.Range("D27").Value = "'159.84"
$ws.Range("E27").Value = "  -0.98%  "

$ws.Range("D28").Value = "'18.43"
$ws.Range("E28").Value = "  -0.39%  "

$ws.Range("D29").Value = "'1.542"
$ws.Range("E29").Value = "  -0.63%  "

$ws.Range("D30").Value = "'4.256"
$ws.Range("E30").Value = "  +0.00%  "

$ws.Range("D31").Value = "'4.201"
$ws.Range("E31").Value = "  -0.14%  "

$ws.Range("E32").Value = "  +0.86%  "

$ws.Range("D33").Value = "'0.05154"
$ws.Range("E33").Value = "  -1.51%  "

$ws.Range("D34").Value = "'0.7947"
$ws.Range("E34").Value = "  +4.47%  "

$ws.Range("D35").Value = "'1.935"
$ws.Range("E35").Value = "  +4.77%  "

$ws.Range("D36").Value = "'1.156"
$ws.Range("E36").Value = "  -1.11%  "

$ws.Range("D37").Value = "'2.692"
$ws.Range("E37").Value = "  -0.66%  "

$ws.Range("D38").Value = "1.336.65"
$ws.Range("E38").Value = "  +9.15%  "

$ws.Range("D39").Value = "'0.01877"
$ws.Range("E39").Value = "  +0.87%  "

$ws.Range("D40").Value = "'2.728"
$ws.Range("E40").Value = "  -0.09%  "

$ws.Range("D41").Value = "'0.9573"
$ws.Range("E41").Value = "  +6.42%  "

$ws.Range("D42").Value = "'6.061"
$ws.Range("E42").Value = "  +9.77%  "

$ws.Range("D43").Value = "'107.79"
$ws.Range("E43").Value = "  -1.06%  "

$ws.Range("D44").Value = "'0.9981"
$ws.Range("E44").Value = "  -0.13%  "

$ws.Range("D45").Value = "'9.797"
$ws.Range("E45").Value = "  +2.61%  "

$ws.Range("D46").Value = "1.995.65"
$ws.Range("E46").Value = "  -0.52%  "

$ws.Range("D47").Value = "'65.09"
$ws.Range("E47").Value = "  -0.38%  "

$ws.Range("D49").Value = "'1.778"
$ws.Range("E49").Value = "  +1.76%  "

$ws.Range("D50").Value = "'0.00000000118"
$ws.Range("E50").Value = "  -5.21%  "

$ws.Range("D51").Value = "'7.023"
$ws.Range("E51").Value = "  -0.29%  "
